$p = $ppt.ActivePresentation

# Slide 10 ("Brief History of the Microsoft Web Stack"):
#   - the "2021" timeline callout becomes "2024"
#   - the "ASP.NET 6" label becomes "ASP.NET 9"
$s = $p.Slides.Item(10)
$s.Shapes.Item(5).TextFrame.TextRange.Text = "2024"
$s.Shapes.Item(13).TextFrame.TextRange.Text = "ASP.NET 9"
